# Actualizacion desde MV -datos-
# Appends six new daily exchange-rate rows (10, 13, 14, 15, 16 and
# 20 September 2021) below the existing data table in Sheet1, extending
# the used range from A1:BH175 to A1:BH181.
#
# Column A holds the date label as literal text (e.g. "10-09-2021"),
# matching the existing rows above it, and columns B:BH hold the numeric
# rate values for each currency series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per new date: first element is the text date label (column A),
# the rest are the numeric values for columns B..BH in order.
$newRows = @(
    @("10-09-2021", 24.21, 791.28, 0, 115.35, 1.26, 36.81, 125.91, 6.21, 91.2, 92.06999999999999, 1127.34, 215.49, 89.81, 583.24, 625.86, 791.28, 791.28, 379.66, 964.98, 791.28, 589.45, 101.75, 562.75, 28.59, 0.03, 936.2, 2.67, 7.84, 863.1799999999999, 0.11, 29.74, 189.47, 50.35, 1095.65, 93.81999999999999, 193.07, 8.07, 0.21, 32.97, 13.99, 15.85, 39.68, 18.59, 102.72, 55.76, 149.3, 0.02, 211.01, 190.81, 10.85, 10.77, 0.06, 4.72, 246.95, 1.86, 0.68, 7.22, 122.66, 206.27),
    @("13-09-2021", 24.16, 789.91, 0, 115.15, 1.26, 36.82, 125.51, 6.18, 91.29000000000001, 91.61, 1126.19, 215.12, 88.45999999999999, 582.01, 624.48, 789.91, 789.91, 379.31, 963.3, 789.91, 589.05, 101.55, 562.73, 28.61, 0.03, 933.15, 2.67, 7.83, 860.37, 0.11, 29.58, 188.91, 50.41, 1093.91, 93.45999999999999, 192.74, 8.050000000000001, 0.21, 32.91, 13.98, 15.85, 39.77, 18.53, 102.5, 55.75, 150.36, 0.02, 210.64, 191.08, 10.8, 10.75, 0.06, 4.71, 246.88, 1.86, 0.68, 7.19, 122.65, 205.53),
    @("14-09-2021", 23.87, 784.26, 0, 114.32, 1.26, 36.51, 124.51, 6.13, 90.73999999999999, 91.05, 1115.59, 213.58, 87.56, 577.3, 619.1900000000001, 784.26, 784.26, 376.36, 956.41, 784.26, 584.1799999999999, 100.81, 557.28, 28.34, 0.03, 925.8200000000001, 2.64, 7.76, 850.52, 0.11, 29.56, 187.15, 50.08, 1084.73, 93, 191.2, 7.99, 0.2, 32.68, 13.88, 15.7, 39.47, 18.37, 101.81, 55.39, 150.04, 0.02, 209.14, 189.02, 10.78, 10.64, 0.06, 4.67, 244.42, 1.84, 0.67, 7.13, 121.71, 203.4),
    @("15-09-2021", 23.79, 783.25, 0, 114.18, 1.25, 36.54, 124.53, 6.13, 90.79000000000001, 91.14, 1115.42, 213.3, 88.8, 574.0599999999999, 618.39, 783.25, 783.25, 375.8, 955.1799999999999, 783.25, 583.6, 100.7, 556.05, 28.3, 0.03, 925.9400000000001, 2.65, 7.76, 852.84, 0.11, 29.55, 187.21, 50, 1083.93, 92.77, 190.85, 7.98, 0.2, 32.64, 13.85, 15.69, 39.36, 18.4, 101.68, 54.79, 149.48, 0.02, 208.87, 188.21, 10.74, 10.63, 0.05, 4.67, 244.22, 1.84, 0.67, 7.15, 121.75, 203.72),
    @("16-09-2021", 23.77, 781.85, 0, 113.97, 1.25, 36.47, 124.16, 6.17, 91.12, 91.06999999999999, 1114.86, 212.92, 87.44, 572.95, 618.36, 781.85, 781.85, 374.34, 953.48, 781.85, 583.08, 100.47, 555.41, 28.26, 0.03, 923.1900000000001, 2.64, 7.75, 849.9299999999999, 0.11, 29.36, 186.64, 49.95, 1081.4, 92.73, 190.39, 7.96, 0.2, 32.58, 13.85, 15.71, 39.34, 18.4, 101.26, 54.27, 148.74, 0.02, 208.49, 188.08, 10.8, 10.64, 0.05, 4.62, 244.12, 1.84, 0.67, 7.15, 121.66, 202.3),
    @("20-09-2021", 23.55, 780.59, 0, 113.79, 1.25, 36.25, 123.42, 6.08, 90.3, 90.38, 1110.84, 212.58, 87.05, 568.65, 615.17, 780.59, 780.59, 374.22, 951.9400000000001, 780.59, 580.15, 100.29, 551.8099999999999, 28.18, 0.03, 917.6900000000001, 2.61, 7.72, 842.51, 0.11, 29.27, 185.54, 49.85, 1075.49, 91.61, 189.79, 7.94, 0.2, 32.52, 13.84, 15.63, 39.12, 18.32, 101.05, 53.44, 148.38, 0.02, 208.16, 187.78, 10.76, 10.62, 0.05, 4.66, 243, 1.84, 0.67, 7.12, 120.97, 200.63)
)

$firstNewRow = 176

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $firstNewRow + $i
    $rowValues = $newRows[$i]

    # Column A: write the date label as text. Assigning the bare string
    # would make Excel auto-detect it as a date and store a date serial,
    # so the cell is briefly forced to the "Text" number format while the
    # value is entered, then restored to the workbook's default style so
    # the cell ends up unstyled - exactly like the rest of column A.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowValues[0]
    $cellA.Style = "Normal"

    # Columns B..BH: plain numeric values.
    for ($j = 1; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, 1 + $j).Value = $rowValues[$j]
    }
}
